$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G7").Value = "2016-08-15 22:39:30"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H7").Value = "2016-08-15 22:39:25"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("H7").Value = "2016-08-15 22:39:30"
